$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-01 Friday" "2025-08-02 Saturday"
Replace-Text "459÷3=" "296÷7="
Replace-Text "203÷7=" "261÷4="
Replace-Text "637÷4=" "331÷7="
Replace-Text "568÷3=" "946÷6="
Replace-Text "596÷5=" "438÷8="
Replace-Text "822÷9=" "684÷3="
Replace-Text "533÷7=" "137÷4="
Replace-Text "125÷8=" "911÷8="
Replace-Text "372÷5=" "346÷5="
Replace-Text "726÷3=" "919÷2="
Replace-Text "600÷2=" "954÷9="
Replace-Text "990÷4=" "970÷4="
Replace-Text "861÷4=" "135÷2="
Replace-Text "656÷2=" "719÷9="
Replace-Text "950÷6=" "320÷8="
Replace-Text "726÷9=" "477÷7="
Replace-Text "179÷5=" "962÷6="
Replace-Text "830÷5=" "777÷8="
Replace-Text "344÷5=" "956÷2="
Replace-Text "439÷8=" "625÷6="
Replace-Text "146÷5=" "769÷3="
Replace-Text "365÷6=" "486÷4="
Replace-Text "182÷7=" "752÷7="
Replace-Text "512÷7=" "507÷9="
Replace-Text "723÷9=" "176÷9="
